$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated statistics (regenerated to filter save games)
$values = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243; E = 6.48142807727062;   G = 28.30127388105354 }
    3 = @{ B = 0.06328177979961902; C = 0.3375848360084654; D = 3.082599426703578; E = 0.4998867070740569; G = 3.98335274958572 }
    4 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 3.082599426703578; E = 0.4998867070740569; G = 6.741336633845642 }
    5 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
